# Apply the pipeline table updates described by the diff:
#   F2: last_run timestamp updated
#   B4: pipeline name changed (SIMCE Lenguaje (Copia) -> Proceso de prueba)
#   C4: description changed (Workflow SIMCE -> Testing pipeline)
#   E4: output format changed (EXCEL -> PDF)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = "2026-01-30 20:30:28"

$ws.Range("B4").Value = "Proceso de prueba"
$ws.Range("C4").Value = "Testing pipeline"
$ws.Range("E4").Value = "PDF"
